$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column P (rows 3-10) to column Q (new 2023 column)
$ws.Range("P3:P10").Copy()
$ws.Range("Q3:Q10").PasteSpecial(-4122)  # xlPasteFormats

# Set header value for 2023
$ws.Range("Q4").Value = 2023

# Set data values for the new column
$ws.Range("Q6").Value = 1209
$ws.Range("Q7").Value = "-"
$ws.Range("Q8").Value = 373
$ws.Range("Q9").Value = 115
$ws.Range("Q10").Value = 781

# Adjust row 5 height
$ws.Rows(5).RowHeight = 27

# Reset selection to A1 (clears the stale P7 selection marker)
$ws.Range("A1").Select() | Out-Null
